# FIX: review text fixes
# Applies the textual corrections from the commit "FIX: review text fixes"
# to the review (рецензия) document.

$d = $word.ActiveDocument

# 1) Title line: "«Разработка сетевых компонентов и их интеграция в шаблон
#    многопользовательской игры на" -> "«Разработка и интеграция сетевых
#    компонентов в шаблон многопользовательской игры на"
$d.Content.Find.Execute(
    "«Разработка сетевых компонентов и их интеграция в шаблон многопользовательской игры на",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "«Разработка и интеграция сетевых компонентов в шаблон многопользовательской игры на",
    2) | Out-Null

# 2) Body sentence: "сетевых компонентов, а также их интеграции в шаблон
#    многопользовательской" -> "и интеграции сетевых компонентов в шаблон
#    многопользовательской"
$d.Content.Find.Execute(
    "сетевых компонентов, а также их интеграции в шаблон многопользовательской",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "и интеграции сетевых компонентов в шаблон многопользовательской",
    2) | Out-Null

# 3) "...в большинстве подобных проектах." -> "...в большинстве подобных проектов."
$d.Content.Find.Execute(
    "большинстве подобных проектах.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "большинстве подобных проектов.",
    2) | Out-Null

# 4) Drop the trailing remark about class descriptions:
#    "...посредством RPC и наличия описания используемых классов."
#    -> "...посредством RPC."
$d.Content.Find.Execute(
    "RPC и наличия описания используемых классов.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RPC.",
    2) | Out-Null

# 5) Remove the stray "_GoBack" bookmark left over from the last edit session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 6) Clean up the spell-check markers around "Техномаш" / "к.т.н" by
#    rewriting that whole line; this merges the runs Word had split for
#    its proofing marks back into plain text runs.
$d.Content.Find.Execute(
    "НПО «Техномаш им. С.А. Афанасьева», к.т.н",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "НПО «Техномаш им. С.А. Афанасьева», к.т.н",
    2) | Out-Null
